$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Add two new "Title and Content" slides at the end of the deck
# (layout 2 == the "Title and Content" custom layout used throughout
# this deck, e.g. the preceding slide 19).
# ---------------------------------------------------------------------

$s1 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s2 = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Slide "Results / Lessons learned 1" -------------------------------
$tr = $s1.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Results"
$tr.LanguageID = "de-DE"
$r = $tr.InsertAfter(" / ")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter("Lessons")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter(" ")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter("learned")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter(" 1")
$r.LanguageID = "de-DE"

# --- Slide "Results / Lessons learned 2" -------------------------------
$tr = $s2.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Results"
$tr.LanguageID = "de-DE"
$r = $tr.InsertAfter(" / ")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter("Lessons")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter(" ")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter("learned")
$r.LanguageID = "de-DE"
$r = $r.InsertAfter(" 2")
$r.LanguageID = "de-DE"
